$d = $word.ActiveDocument
$s = $d.Styles.Item("SubtitleChar")
Write-Host "Font:" $s.Font
Write-Host "Font.Color:" $s.Font.Color
Write-Host "Font.ColorIndex:" $s.Font.ColorIndex
Write-Host "Font.ThemeColor:" $s.Font.ThemeColor
